$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Update the two vehicle-type data sheets (BCbVT-passenger, BCbVT-freight)
# to add two new vehicle-type columns (LPG vehicle, hydrogen vehicle),
# move the existing "nonroad vehicle" column out to make room, and label
# the corner cell with the battery-capacity units header.
# ---------------------------------------------------------------------
$sheetNames = @("BCbVT-passenger", "BCbVT-freight")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # The last existing header (column G) was "nonroad vehicle". Shift
    # it out to column I, and put the two new vehicle types into the
    # freed-up G and H columns, matching the header style used by the
    # rest of the row (bold, right-aligned).
    $g1 = $ws.Cells.Item(1, 7)
    $g1.Value = "LPG vehicle"
    $g1.Font.Bold = $true
    $g1.HorizontalAlignment = -4152

    $h1 = $ws.Cells.Item(1, 8)
    $h1.Value = "hydrogen vehicle"
    $h1.Font.Bold = $true
    $h1.HorizontalAlignment = -4152

    $i1 = $ws.Cells.Item(1, 9)
    $i1.Value = "nonroad vehicle"
    $i1.Font.Bold = $true
    $i1.HorizontalAlignment = -4152

    # Corner cell A1: add the row-header label and keep the existing
    # right-aligned style, but also turn word-wrap on for it.
    $a1 = $ws.Cells.Item(1, 1)
    $a1.Value = "Battery Capacity (MW*hr/vehicle"
    $a1.WrapText = $true
    $a1.HorizontalAlignment = -4152

    # Fill the two new columns with zero values for every data row,
    # consistent with the other vehicle-type columns.
    for ($r = 2; $r -le 7; $r++) {
        $ws.Cells.Item($r, 8).Value = 0
        $ws.Cells.Item($r, 9).Value = 0
    }
}
